# "JS in methods for Demo"
# Refresh the first 5 doctor records on the "Doctors" sheet with a new
# batch of Gynecologist/Obstetrician doctors (Delhi locations).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

# Row 2: Dr. Karishma Bhatia
$ws.Range("A2").Value = "Dr. Karishma Bhatia"
$ws.Range("B2").Value = "Gynecologist/Obstetrician"
$ws.Range("C2").Value = "13 years experience overall"
$ws.Range("D2").Value = "Rohini,Delhi"

# Row 3: Dr. Preeti Tahilyani
$ws.Range("A3").Value = "Dr. Preeti Tahilyani"
$ws.Range("B3").Value = "Gynecologist/Obstetrician"
$ws.Range("C3").Value = "22 years experience overall"
$ws.Range("D3").Value = "Dwarka,Delhi"

# Row 4: Dr. Tejashri Shrotri
$ws.Range("A4").Value = "Dr. Tejashri Shrotri"
$ws.Range("B4").Value = "Gynecologist/Obstetrician"
$ws.Range("C4").Value = "20 years experience overall"
$ws.Range("D4").Value = "Hauz Khas,Delhi"

# Row 5: Dr. Shalini Chaudhry
$ws.Range("A5").Value = "Dr. Shalini Chaudhry"
$ws.Range("B5").Value = "Gynecologist/Obstetrician"
$ws.Range("C5").Value = "27 years experience overall"
$ws.Range("D5").Value = "Dwarka,Delhi"

# Row 6: Dr. Shakuntla Shukla
$ws.Range("A6").Value = "Dr. Shakuntla Shukla"
$ws.Range("B6").Value = "Gynecologist/Obstetrician"
$ws.Range("C6").Value = "41 years experience overall"
$ws.Range("D6").Value = "Lajpat Nagar,Delhi"
